$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: add "Completed?" = "Yes?" (C7, styled like existing C4/C9 "Yes" cells)
#        and "Completion Date" = "30/05/2024 (TEST)" as plain text (D7, styled like D4)
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Range("C7").Value = "Yes?"

$ws.Range("D4").Copy() | Out-Null
$ws.Range("D7").PasteSpecial(-4122) | Out-Null
$ws.Range("D7").Value = "30/05/2024 (TEST)"

# Row 8: add "Completed?" = "Yes" (C8) and "Completion Date" = 30/05/2024 as a real date (D8)
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C8").PasteSpecial(-4122) | Out-Null
$ws.Range("C8").Value = "Yes"

$ws.Range("D9").Copy() | Out-Null
$ws.Range("D8").PasteSpecial(-4122) | Out-Null
$ws.Range("D8").Value = 45442

# Row 13: add "Completed?" = "Yes" (C13) and "Completion Date" = 29/05/2024 as a real date (D13)
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Value = "Yes"

$ws.Range("D9").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("D13").Value = 45441

# Row 20: rename the task text, and add "Completed?" = "Yes" (C20) and
#         "Completion Date" = 30/05/2024 as a real date (D20)
$ws.Range("A20").Value = "Add map and fps meter to cl_showinfo"

$ws.Range("C9").Copy() | Out-Null
$ws.Range("C20").PasteSpecial(-4122) | Out-Null
$ws.Range("C20").Value = "Yes"

$ws.Range("D9").Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4122) | Out-Null
$ws.Range("D20").Value = 45442

# Update the active selection to match the final workbook state (A7)
$ws.Range("A7").Select() | Out-Null

$excel.CutCopyMode = 0
